$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Resolving-Mac -----------------------------------------
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 25.28841266666667
$ws.Range("H2").Value = 75.86523800000001
$ws.Range("I2").Value = 0.08258585054448338
$ws.Range("J2").Value = 0.08258585054448338
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009727666666666667
$ws.Range("N2").Value = 0.029183
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.2459972489504445
$ws.Range("R2").Value = 2.213975240554
$ws.Range("S2").Value = 0.08258585054448338
$ws.Range("T2").Value = 0.08258585054448338

# --- Row 3: FAPs -> Resolving-Mac -----------------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 195.050573
$ws.Range("H3").Value = 585.1517190000001
$ws.Range("I3").Value = 0.6369880815661784
$ws.Range("J3").Value = 0.6369880815661784
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.009727666666666667
$ws.Range("N3").Value = 0.029183
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.897386957286334
$ws.Range("R3").Value = 17.076482615577
$ws.Range("S3").Value = 0.6369880815661784
$ws.Range("T3").Value = 0.6369880815661784

# --- Row 4: MuSCs -> Resolving-Mac -----------------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 55.14511
$ws.Range("H4").Value = 165.43533
$ws.Range("I4").Value = 0.180090615917626
$ws.Range("J4").Value = 0.180090615917626
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009727666666666667
$ws.Range("N4").Value = 0.029183
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.5364332483766666
$ws.Range("R4").Value = 4.82789923539
$ws.Range("S4").Value = 0.180090615917626
$ws.Range("T4").Value = 0.180090615917626

# --- Row 5: Resolving-Mac -> Resolving-Mac ---------------------------------
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 30.723475
$ws.Range("H5").Value = 92.17042499999999
$ws.Range("I5").Value = 0.1003354519717122
$ws.Range("J5").Value = 0.1003354519717122
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009727666666666667
$ws.Range("N5").Value = 0.029183
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.2988677236416666
$ws.Range("R5").Value = 2.689809512775
$ws.Range("S5").Value = 0.1003354519717122
$ws.Range("T5").Value = 0.1003354519717122

# --- Remove the now-obsolete rows 6-9 (MuSCs/Resolving-Mac x ECs pairs) ----
$ws.Range("A6:T9").EntireRow.Delete()
